# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 5
    4  = 7
    5  = 7
    6  = 6
    7  = 8
    8  = 6
    9  = 5
    10 = 6
    11 = 4
    12 = 7
    13 = 4
    14 = 6
    15 = 4
    16 = 6
    17 = 6
    18 = 3
    19 = 4
    20 = 4
    21 = 9
    22 = 6
    23 = 8
    24 = 7
    25 = 5
    26 = 5
    27 = 3
    28 = 6
    29 = 4
    30 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
